$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 7
$lastRow = 36

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Extend existing Week 1-4 formulas (columns L:O) from $AT$32 to $BH$32
    foreach ($col in @("L", "M", "N", "O")) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        $cell.Formula = $f.Replace("`$AT`$32", "`$BH`$32")
    }

    # Match formatting of P/Q to the neighboring calculated column (O) before
    # filling in their formulas, mirroring how Excel extends a table's
    # calculated-column formatting when new calculated cells are created.
    $ws.Range("O" + $r).Copy()
    $ws.Range("P" + $r + ":Q" + $r).PasteSpecial(-4122)

    # Add new formulas for Week 5 (P) and Week 6 (Q)
    $ws.Range("P" + $r).Formula = "=SUM(INDEX(`$R`$1:`$BH`$32,MATCH(`$A" + $r + ",`$R`$1:`$R`$32,0),30):INDEX(`$R`$1:`$BH`$32,MATCH(`$A" + $r + ",`$R`$1:`$R`$32,0),36))"
    $ws.Range("Q" + $r).Formula = "=SUM(INDEX(`$R`$1:`$BH`$32,MATCH(`$A" + $r + ",`$R`$1:`$R`$32,0),37):INDEX(`$R`$1:`$BH`$32,MATCH(`$A" + $r + ",`$R`$1:`$R`$32,0),43))"
}

$excel.CutCopyMode = $false

# Update the active selection
$ws.Range("M31").Select()
